$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the missing "pair_kind" value ("generic") for the practice rows (rows 2-5)
$ws.Range("J2:J5").Value = "generic"

# --- New "stim details" block added at the bottom of the sheet ---
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

$data = @(
    @(6, "video"),
    @(6, "video"),
    @(7, "video"),
    @(7, "video"),
    @(6, "audio"),
    @(6, "audio"),
    @(7, "audio"),
    @(7, "audio")
)

$r = 29
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r++
}
